# Apply trade #78 closing update to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 78      # Total Trades
$wsSummary.Range("B9").Value = 32.05   # Win Rate %

# --- Strategy Status sheet -------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 78       # MarketMaking Total Trades
$wsStatus.Range("G4").Value = 32.05    # MarketMaking Win Rate %

# --- Append the new trade row (#78) to "All Trades" and "MarketMaking" ----
$newRow = 79

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = 78

    # Force the date column to stay a plain text string instead of being
    # auto-converted to a date serial number by Excel's COM layer.
    $dateCell = $ws.Cells.Item($newRow, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()

    $ws.Cells.Item($newRow, 3).Value = "15:49:28"
    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"
    $ws.Cells.Item($newRow, 5).Value = "UP"
    $ws.Cells.Item($newRow, 6).Value = 0.92
    $ws.Cells.Item($newRow, 7).Value = 0.92
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"
    $ws.Cells.Item($newRow, 9).Value = 0
    $ws.Cells.Item($newRow, 10).Value = 0
    $ws.Cells.Item($newRow, 11).Value = 99.90000000000001
    $ws.Cells.Item($newRow, 12).Value = 0
    $ws.Cells.Item($newRow, 13).Value = 0
    $ws.Cells.Item($newRow, 14).Value = 0.6
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($newRow, 16).Value = "early_exit"
    $ws.Cells.Item($newRow, 17).Value = 0.11
}
